$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.118.89'
$ws.Range("E2").Value = '  -1.53%  '

$ws.Range("D3").Value = '1.780.58'
$ws.Range("E3").Value = '  -1.89%  '

$ws.Range("E4").Value = '  +0.39%  '

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '335.99'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  -2.47%  '

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  +0.30%  '

$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3835'
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = '  -0.02%  '

$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3415'
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = '  -3.03%  '

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.17'
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = '  -3.09%  '

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.190'
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = '  -3.64%  '

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07465'
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = '  -4.09%  '

$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.004'
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = '  +0.30%  '

$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.65'
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = '  -3.42%  '

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.417'
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = '  -3.04%  '

$ws.Range("D15").Value = '1.780.63'
$ws.Range("E15").Value = '  -1.88%  '

$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.068'
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = '  -2.12%  '

$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001086'
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = '  -4.01%  '

$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06653'
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = '  -1.03%  '

$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '83.35'
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = '  -3.67%  '

$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.003'
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = '  +0.32%  '

$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.562'
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = '  +0.40%  '

$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.27'
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = '  -2.93%  '

$ws.Range("D23").Value = '27.120.67'
$ws.Range("E23").Value = '  -1.48%  '

$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.24'
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = '  -7.08%  '

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.386'
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = '  -3.18%  '

$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.515'
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = '  -6.34%  '

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.469'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  -2.06%  '

$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.15'
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = '  -3.96%  '

$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '154.16'
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = '  +0.19%  '

$ws.Range("D30").Value = '1.982.60'
$ws.Range("E30").Value = '  -1.76%  '

$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '134.02'
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = '  -1.89%  '

$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.021'
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = '  -1.55%  '

$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.016'
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  -5.94%  '

$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08677'
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = '  -1.64%  '

$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.15'
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = '  -6.40%  '

$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.631'
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = '  -4.95%  '

$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6845'
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = '  -3.68%  '

$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.388'
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = '  -4.70%  '

$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06292'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  -4.02%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.731'
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = '  -3.44%  '

$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2177'
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = '  -4.39%  '

$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.02320'
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = '  -4.56%  '

$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.233'
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = '  -4.21%  '

$ws.Range("E44").Value = '  -4.56%  '

$ws.Range("E45").Value = '  +0.31%  '

$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6424'
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = '  -3.45%  '

$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.848'
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = '  -2.98%  '

$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.130'
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = '  -2.59%  '

$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '130.39'
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = '  -2.23%  '

$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07109'
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = '  -3.42%  '

$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.61'
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = '  -2.72%  '
